$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.7287194209349384;  C = 0.05231270169004087;  D = 0.7127328510149897; E = 0.4998867070740569; G = 1.993651680714026 }
    3 = @{ B = 3.182878228561681;   C = 9.226618575922256;    D = 0.1529057820181812; E = 6.48142807727062;    G = 19.04383066377274 }
    4 = @{ B = 3.182878228561681;   C = 1.65323645889881;     D = 3.082599426703578;  E = 6.48142807727062;    G = 14.40014219143469 }
    5 = @{ B = 0.06328177979961902; C = 0.3375848360084654;   D = 0.1529057820181812; E = 0.4998867070740569; G = 1.053659104900323 }
    6 = @{ B = 3.182878228561681;   C = 1.65323645889881;     D = 3.082599426703578;  E = 6.48142807727062;    G = 14.40014219143469 }
    7 = @{ B = 0.06328177979961902; C = 0.0001537489499301437; D = 0.7127328510149897; E = 6.48142807727062;   G = 7.257596457035159 }
    8 = @{ B = 0.7287194209349384;  C = 1.65323645889881;     D = 0.1529057820181812; E = 0.4998867070740569; G = 3.034748368925986 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
